$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins, Losses, Ties), reusing the same
# formatting (bold, centered, bordered) as the rest of the header row by
# copying an existing header cell's format before overwriting its text.
$ws.Range("AB1").Copy($ws.Range("AC1"))
$ws.Range("AB1").Copy($ws.Range("AD1"))
$ws.Range("AB1").Copy($ws.Range("AE1"))

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the team record (Wins=78, Losses=84, Ties=0) for every player row.
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 78
    $ws.Cells.Item($r, 30).Value = 84
    $ws.Cells.Item($r, 31).Value = 0
}
